$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 17
$ws.Range("H17").Value = 3576.5405
$ws.Range("J17").Value = 3824.4707
$ws.Range("L17").Value = 11473.4121
$ws.Range("N17").Value = -11809.4121

# Row 33
$ws.Range("H33").Value = 469.8
$ws.Range("I33").Value = 577.5
$ws.Range("J33").Value = 398
$ws.Range("K33").Value = 577.5
$ws.Range("L33").Value = 398
$ws.Range("M33").Value = -348.5
$ws.Range("N33").Value = -856

# Row 34
$ws.Range("H34").Value = 7849.375
$ws.Range("J34").Value = 10998.75
$ws.Range("L34").Value = 10998.75
$ws.Range("N34").Value = -11404.75

# Row 36
$ws.Range("H36").Value = 7849.375
$ws.Range("J36").Value = 10998.75
$ws.Range("L36").Value = 10998.75
$ws.Range("N36").Value = -12428.75

# Row 69
$ws.Range("H69").Value = 8074.8945
$ws.Range("I69").Value = 2004.3334
$ws.Range("K69").Value = 6013.0002
$ws.Range("M69").Value = -5139.0002

# Row 72
$ws.Range("H72").Value = 8074.8945
$ws.Range("I72").Value = 2004.3334
$ws.Range("K72").Value = 18039.0006
$ws.Range("M72").Value = -13671.0006

# Row 138
$ws.Range("H138").Value = 3479.8708
$ws.Range("I138").Value = 4198.3335
$ws.Range("J138").Value = 3307.44
$ws.Range("K138").Value = 12595.0005
$ws.Range("L138").Value = 9922.32
$ws.Range("M138").Value = -7455.000499999998
$ws.Range("N138").Value = -20202.32

# Row 141
$ws.Range("H141").Value = 1068.5
$ws.Range("I141").Value = 1079
$ws.Range("J141").Value = 995
$ws.Range("K141").Value = 3237
$ws.Range("L141").Value = 2985
$ws.Range("M141").Value = 1943
$ws.Range("N141").Value = -13345

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Range("H61").Value = 6292.8125
$ws.Range("I61").Value = 5857.0967
$ws.Range("K61").Value = 5857.0967
$ws.Range("M61").Value = -5645.0967

# Row 74
$ws.Range("H74").Value = 18520392
$ws.Range("I74").Value = 23811572
$ws.Range("J74").Value = 1256.75
$ws.Range("K74").Value = 23811572
$ws.Range("L74").Value = 1256.75
$ws.Range("M74").Value = -23810698
$ws.Range("N74").Value = -3004.75

# Row 77
$ws.Range("H77").Value = 18520392
$ws.Range("I77").Value = 23811572
$ws.Range("J77").Value = 1256.75
$ws.Range("K77").Value = 119057860
$ws.Range("L77").Value = 6283.75
$ws.Range("M77").Value = -119053492
$ws.Range("N77").Value = -15019.75

# Row 132
$ws.Range("H132").Value = 2915.818
$ws.Range("I132").Value = 1977.9231
$ws.Range("K132").Value = 5933.7693
$ws.Range("M132").Value = -3403.7693

# Row 136
$ws.Range("H136").Value = 6292.8125
$ws.Range("I136").Value = 5857.0967
$ws.Range("K136").Value = 17571.2901
$ws.Range("M136").Value = -15021.2901

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 99
$ws.Range("H99").Value = 1991.5
$ws.Range("I99").Value = 2399.5
$ws.Range("J99").Value = 1787.5
$ws.Range("K99").Value = 2399.5
$ws.Range("L99").Value = 1787.5
$ws.Range("M99").Value = -901.5
$ws.Range("N99").Value = -4783.5

# Row 105
$ws.Range("H105").Value = 20043.611
$ws.Range("J105").Value = 20083.166
$ws.Range("L105").Value = 20083.166
$ws.Range("N105").Value = -23577.166

# Row 134
$ws.Range("H134").Value = 1863.375
$ws.Range("I134").Value = 1291.6086
$ws.Range("K134").Value = 3874.8258
$ws.Range("M134").Value = -1339.8258

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 19
$ws.Range("H19").Value = 1976
$ws.Range("I19").Value = 1284.5
$ws.Range("J19").Value = 2227.4546
$ws.Range("K19").Value = 1284.5
$ws.Range("L19").Value = 2227.4546
$ws.Range("M19").Value = -1114.5
$ws.Range("N19").Value = -2567.4546

# Row 24
$ws.Range("H24").Value = 1976
$ws.Range("I24").Value = 1284.5
$ws.Range("J24").Value = 2227.4546
$ws.Range("K24").Value = 1284.5
$ws.Range("L24").Value = 2227.4546
$ws.Range("M24").Value = -1114.5
$ws.Range("N24").Value = -2567.4546

# Row 31
$ws.Range("H31").Value = 27849.273
$ws.Range("J31").Value = 109865.5
$ws.Range("L31").Value = 109865.5
$ws.Range("N31").Value = -110455.5

# Row 34
$ws.Range("H34").Value = 27849.273
$ws.Range("J34").Value = 109865.5
$ws.Range("L34").Value = 109865.5
$ws.Range("N34").Value = -110269.5

# Row 107
$ws.Range("H107").Value = 1555.9333
$ws.Range("I107").Value = 2149.1667
$ws.Range("J107").Value = 1160.4445
$ws.Range("K107").Value = 2149.1667
$ws.Range("L107").Value = 1160.4445
$ws.Range("M107").Value = -229.1667000000002
$ws.Range("N107").Value = -5000.4445

# Row 132
$ws.Range("H132").Value = 1269.1154
$ws.Range("I132").Value = 1269.1154
$ws.Range("K132").Value = 3807.3462
$ws.Range("M132").Value = -1277.3462

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 86
$ws.Range("H86").Value = 2733.111
$ws.Range("J86").Value = 3866
$ws.Range("L86").Value = 11598
$ws.Range("N86").Value = -13970

# Row 89
$ws.Range("H89").Value = 2733.111
$ws.Range("J89").Value = 3866
$ws.Range("L89").Value = 34794
$ws.Range("N89").Value = -46650

# Row 107
$ws.Range("H107").Value = 1246.6875
$ws.Range("I107").Value = 1730
$ws.Range("J107").Value = 1057.5652
$ws.Range("K107").Value = 5190
$ws.Range("L107").Value = 3172.6956
$ws.Range("M107").Value = -3270
$ws.Range("N107").Value = -7012.6956

# Row 108
$ws.Range("H108").Value = 7463.75
$ws.Range("I108").Value = 2900
$ws.Range("J108").Value = 12027.5
$ws.Range("K108").Value = 8700
$ws.Range("L108").Value = 36082.5
$ws.Range("M108").Value = -5820
$ws.Range("N108").Value = -41842.5

# Row 119
$ws.Range("H119").Value = 37521.2
$ws.Range("I119").Value = 37673
$ws.Range("J119").Value = 37293.5
$ws.Range("K119").Value = 113019
$ws.Range("L119").Value = 111880.5
$ws.Range("M119").Value = -108181
$ws.Range("N119").Value = -121556.5

# Row 121
$ws.Range("H121").Value = 4383.25
$ws.Range("J121").Value = 9033
$ws.Range("L121").Value = 27099
$ws.Range("N121").Value = -29719

# Row 126
$ws.Range("H126").Value = 4516
$ws.Range("J126").Value = 4516
$ws.Range("L126").Value = 13548
$ws.Range("N126").Value = -23428

# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# Row 131
$ws.Range("H131").Value = 24574422
$ws.Range("I131").Value = 17859408
$ws.Range("J131").Value = 32408608
$ws.Range("K131").Value = 53578224
$ws.Range("L131").Value = 97225824
$ws.Range("M131").Value = -53573184
$ws.Range("N131").Value = -97235904

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 126
$ws.Range("H126").Value = 4111.4546
$ws.Range("I126").Value = 2621.2
$ws.Range("K126").Value = 7863.599999999999
$ws.Range("M126").Value = -5393.599999999999

# Row 132
$ws.Range("H132").Value = 10412.565
$ws.Range("I132").Value = 8832.2
$ws.Range("J132").Value = 13375.75
$ws.Range("K132").Value = 26496.6
$ws.Range("L132").Value = 40127.25
$ws.Range("M132").Value = -23966.6
$ws.Range("N132").Value = -45187.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 46
$ws.Range("H46").Value = 4750
$ws.Range("J46").Value = 5055.5557
$ws.Range("L46").Value = 5055.5557
$ws.Range("N46").Value = -5431.5557

# Row 68
$ws.Range("H68").Value = 4475.2
$ws.Range("I68").Value = 3219
$ws.Range("J68").Value = 5312.6665
$ws.Range("K68").Value = 3219
$ws.Range("L68").Value = 5312.6665
$ws.Range("M68").Value = -2470
$ws.Range("N68").Value = -6810.6665

# Row 71
$ws.Range("H71").Value = 4475.2
$ws.Range("I71").Value = 3219
$ws.Range("J71").Value = 5312.6665
$ws.Range("K71").Value = 16095
$ws.Range("L71").Value = 26563.3325
$ws.Range("M71").Value = -12351
$ws.Range("N71").Value = -34051.3325

# Row 93
$ws.Range("H93").Value = 4249.5
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

# Row 111
$ws.Range("H111").Value = 73684
$ws.Range("J111").Value = 73684
$ws.Range("L111").Value = 73684
$ws.Range("N111").Value = -81864

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 107
$ws.Range("H107").Value = 872.6875
$ws.Range("I107").Value = 685.6667
$ws.Range("J107").Value = 1433.75
$ws.Range("K107").Value = 2057.0001
$ws.Range("L107").Value = 4301.25
$ws.Range("M107").Value = -137.0001000000002
$ws.Range("N107").Value = -8141.25

# Row 122
$ws.Range("H122").Value = 7911.6055
$ws.Range("I122").Value = 3629.4546
$ws.Range("J122").Value = 13799.5625
$ws.Range("K122").Value = 10888.3638
$ws.Range("L122").Value = 41398.6875
$ws.Range("M122").Value = -8438.3638
$ws.Range("N122").Value = -46298.6875
